## TC01_Canine_CaseFiles_Study-COTB.xlsx -> adds a "TabName" column and
## duplicates the query row for the Cases/Samples/Files tabs, refreshes the
## Cypher query text and the Neo4j/Web export file names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shift the existing A:D data one column to the right, making room for
#     the new "TabName" lead column --------------------------------------
$ws.Columns.Item(1).Insert()

# --- row 1 (headers) -----------------------------------------------------
$ws.Range("A1").Value = "TabName"
$ws.Range("B1").Value = "query"
$ws.Range("C1").Value = "StatQuery"
$ws.Range("D1").Value = "dbExcel"
$ws.Range("E1").Value = "WebExcel"

# --- the (updated) long Cypher query strings shared by every data row ----
$caseQuery = @"
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE demo.breed IN ['COTC007B']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID``,
        coalesce(s.clinical_study_designation, '') AS ``Study Code``,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease``,
        coalesce(demo.patient_age_at_enrollment, '') AS Age,
        coalesce(demo.sex, '') AS Sex,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@

$statQuery = @"
MATCH (s:study)
  WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies
  MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies
  MATCH (d:diagnosis)
  WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies
  MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
    WHERE s.clinical_study_designation IN ['COTC007B']
  OPTIONAL MATCH (f:file)-[*]->(c)
  OPTIONAL MATCH (samp:sample)-[*]->(c)
  WITH DISTINCT c AS c, p, s, demo, diag, f, samp
  RETURN count(DISTINCT(f)) as number_of_files ,
             count(DISTINCT(samp)) as number_of_sample ,
             count(DISTINCT(c.case_id)) as number_of_cases ,
             count(DISTINCT(s.clinical_study_designation)) as number_of_study
"@

$neo4jFile = "TC03_Canine_Filter_Study-GLIOMA_Neo4jData.xlsx"
$webFile   = "TC03_Canine_Filter_Study-GLIOMA_WebData.xlsx"

# --- row 2: CasesTab -------------------------------------------------------
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $caseQuery
$ws.Range("C2").Value = $statQuery
$ws.Range("D2").Value = $neo4jFile
$ws.Range("E2").Value = $webFile

# --- duplicate row 2 (with its formatting) into rows 3 and 4 -------------
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(4).Insert()

# --- row 3: SamplesTab / row 4: FilesTab ----------------------------------
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("A4").Value = "FilesTab"

# --- formatting: wrap text on the two query columns for every data row ---
$ws.Range("B2:C4").WrapText = $true

# --- row heights for the (now taller, three-line) data rows --------------
$ws.Rows.Item(2).RowHeight = 244.8
$ws.Rows.Item(3).RowHeight = 244.8
$ws.Rows.Item(4).RowHeight = 244.8

# --- column widths ----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10
$ws.Range("B1:C1").ColumnWidth = 75
$ws.Columns.Item(4).ColumnWidth = 69.33333333333333
$ws.Columns.Item(5).ColumnWidth = 37.166666666666664

# --- zoom the view to 115% (as recorded by the refreshed Excel build) ----
$excel.ActiveWindow.Zoom = 115
